$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 3 values (was UserName1/Password1 -> now T1/P1)
$ws.Range("A3").Value = "T1"
$ws.Range("B3").Value = "P1"

# Append a new data row 4 (T2/P2)
$ws.Range("A4").Value = "T2"
$ws.Range("B4").Value = "P2"

# Move/record the active selection on the newly added cell
$ws.Range("B4").Select() | Out-Null
